$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Q4").Value = 2020
